$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 40
$ws.Range("H40").Value = 1470.1428
$ws.Range("I40").Value = 1458.2
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 1458.2
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -1283.2
$ws.Range("N40").Value = -1850

# ALC row 61
$ws.Range("H61").Value = 245.44444
$ws.Range("I61").Value = 213.625
$ws.Range("J61").Value = 500
$ws.Range("K61").Value = 640.875
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -468.875
$ws.Range("N61").Value = -1844

# ALC row 125
$ws.Range("H125").Value = 2140.889
$ws.Range("I125").Value = 2844
$ws.Range("J125").Value = 1789.3334
$ws.Range("K125").Value = 25596
$ws.Range("L125").Value = 16104.0006
$ws.Range("M125").Value = -23136
$ws.Range("N125").Value = -21024.0006

# ALC row 132
$ws.Range("H132").Value = 20918692
$ws.Range("I132").Value = 22312680
$ws.Range("K132").Value = 66938040
$ws.Range("M132").Value = -66935510

# ALC row 138
$ws.Range("H138").Value = 1879.8
$ws.Range("I138").Value = 918.9804
$ws.Range("J138").Value = 2879.8367
$ws.Range("K138").Value = 2756.9412
$ws.Range("L138").Value = 8639.5101
$ws.Range("M138").Value = 2383.0588
$ws.Range("N138").Value = -18919.5101

# ALC row 139
$ws.Range("H139").Value = 41165.688
$ws.Range("J139").Value = 41165.688
$ws.Range("L139").Value = 41165.688
$ws.Range("N139").Value = -51445.688

# ALC row 141
$ws.Range("H141").Value = 2522.739
$ws.Range("I141").Value = 2166.15
$ws.Range("K141").Value = 6498.450000000001
$ws.Range("M141").Value = -1318.450000000001

$ws = $wb.Worksheets.Item("ARM")
# ARM row 31
$ws.Range("H31").Value = 9512.556
$ws.Range("I31").Value = 9512.556
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 9512.556
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -9218.556
$ws.Range("N31").ClearContents()

# ARM row 32
$ws.Range("H32").Value = 8288.384
$ws.Range("I32").Value = 5229.9404
$ws.Range("J32").Value = 14692
$ws.Range("K32").Value = 5229.9404
$ws.Range("L32").Value = 14692
$ws.Range("M32").Value = -4942.9404
$ws.Range("N32").Value = -15266

# ARM row 110
$ws.Range("H110").Value = 3050.6
$ws.Range("I110").Value = 3998.8333
$ws.Range("J110").Value = 1628.25
$ws.Range("K110").Value = 3998.8333
$ws.Range("L110").Value = 1628.25
$ws.Range("M110").Value = -1953.8333
$ws.Range("N110").Value = -5718.25

# ARM row 122
$ws.Range("H122").Value = 2297.7856
$ws.Range("I122").Value = 1280.75
$ws.Range("K122").Value = 3842.25
$ws.Range("M122").Value = -1392.25

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Range("H31").Value = 3216.9707
$ws.Range("I31").Value = 1095.95
$ws.Range("K31").Value = 1095.95
$ws.Range("M31").Value = -800.95

# CRP row 34
$ws.Range("H34").Value = 3216.9707
$ws.Range("I34").Value = 1095.95
$ws.Range("K34").Value = 1095.95
$ws.Range("M34").Value = -893.95

# CRP row 134
$ws.Range("H134").Value = 4481.25
$ws.Range("I134").Value = 5049.6665
$ws.Range("J134").Value = 3344.4167
$ws.Range("K134").Value = 15148.9995
$ws.Range("L134").Value = 10033.2501
$ws.Range("M134").Value = -12613.9995
$ws.Range("N134").Value = -15103.2501

# CRP row 139
$ws.Range("H139").Value = 101688.57
$ws.Range("J139").Value = 101688.57
$ws.Range("L139").Value = 101688.57
$ws.Range("N139").Value = -111968.57

$ws = $wb.Worksheets.Item("CUL")
# CUL row 3
$ws.Range("H3").Value = 4697.304
$ws.Range("I3").Value = 3216.125
$ws.Range("J3").Value = 8082.857
$ws.Range("K3").Value = 9648.375
$ws.Range("L3").Value = 24248.571
$ws.Range("M3").Value = -9536.375
$ws.Range("N3").Value = -24472.571

# CUL row 113
$ws.Range("H113").Value = 629.96075
$ws.Range("I113").Value = 626.32434
$ws.Range("J113").Value = 639.5714
$ws.Range("K113").Value = 1878.97302
$ws.Range("L113").Value = 1918.7142
$ws.Range("M113").Value = 291.0269800000001
$ws.Range("N113").Value = -6258.7142

# CUL row 129
$ws.Range("H129").Value = 2074.1924
$ws.Range("I129").Value = 3650
$ws.Range("J129").Value = 918.6
$ws.Range("K129").Value = 10950
$ws.Range("L129").Value = 2755.8
$ws.Range("M129").Value = -5950
$ws.Range("N129").Value = -12755.8

# CUL row 131
$ws.Range("H131").Value = 9260243
$ws.Range("I131").Value = 166670160
$ws.Range("J131").Value = 836.1961
$ws.Range("K131").Value = 500010480
$ws.Range("L131").Value = 2508.5883
$ws.Range("M131").Value = -500005440
$ws.Range("N131").Value = -12588.5883

$ws = $wb.Worksheets.Item("GSM")
# GSM row 22
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()

# GSM row 58
$ws.Range("H58").Value = 27000
$ws.Range("I58").Value = 27000
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 27000
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -26723
$ws.Range("N58").ClearContents()

# GSM row 70
$ws.Range("H70").Value = 5959.091
$ws.Range("I70").Value = 5696.3228
$ws.Range("J70").Value = 6585.6924
$ws.Range("K70").Value = 5696.3228
$ws.Range("L70").Value = 6585.6924
$ws.Range("M70").Value = -5426.3228
$ws.Range("N70").Value = -7125.6924

# GSM row 73
$ws.Range("H73").Value = 5959.091
$ws.Range("I73").Value = 5696.3228
$ws.Range("J73").Value = 6585.6924
$ws.Range("K73").Value = 5696.3228
$ws.Range("L73").Value = 6585.6924
$ws.Range("M73").Value = -4760.3228
$ws.Range("N73").Value = -8457.6924

# GSM row 102
$ws.Range("H102").Value = 2802.9355
$ws.Range("I102").Value = 2064.913
$ws.Range("K102").Value = 2064.913
$ws.Range("M102").Value = -442.913

# GSM row 113
$ws.Range("H113").Value = 1509.9
$ws.Range("I113").Value = 1442.7142
$ws.Range("J113").Value = 1666.6666
$ws.Range("K113").Value = 1442.7142
$ws.Range("L113").Value = 1666.6666
$ws.Range("M113").Value = 727.2858000000001
$ws.Range("N113").Value = -6006.6666

# GSM row 136
$ws.Range("H136").Value = 13607.72
$ws.Range("J136").Value = 13607.72
$ws.Range("L136").Value = 40823.16
$ws.Range("N136").Value = -45923.16

$ws = $wb.Worksheets.Item("WVR")
# WVR row 133
$ws.Range("H133").Value = 52000
$ws.Range("J133").Value = 52000
$ws.Range("L133").Value = 52000
$ws.Range("N133").Value = -62120

# WVR row 136
$ws.Range("H136").Value = 3463.5881
$ws.Range("I136").Value = 1060.2222
$ws.Range("J136").Value = 6167.375
$ws.Range("K136").Value = 3180.6666
$ws.Range("L136").Value = 18502.125
$ws.Range("M136").Value = -630.6665999999996
$ws.Range("N136").Value = -23602.125
